$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fragen")

# Row 10 (Statistik_Code 33111, Waldfläche): fix question wording from
# "In welcher Stadt..." to "Welche Stadt..."
$ws.Range("G10").Value = "Welche Stadt oder welcher Landkreis hat die meiste Waldfläche?"

# Row 14 (Statistik_Code 46251, Kraftfahrzeugbestand): fix "Merkmal" from
# the generic "insgesamt" to "Personenkraftwagen"
$ws.Range("D14").Value = "Personenkraftwagen"

# Update the active selection to the last-edited cell
$ws.Activate()
$ws.Range("D14").Select()
